{"js": "// Apply literal text replacements for each unique cell/title text.\n// Using Range.search per old text (all texts are unique in this document),\n// then replacing the matched range's text in place so formatting (rPr) is preserved.\nconst replacements = [\n  { oldText: \"2024-11-28 Thursday\", newText: \"2024-11-29 Friday\" },\n  { oldText: \"598\u00f73=199, 1\", newText: \"947\u00f79=105, 2\" },\n  { oldText: \"612\u00f74=153, 0\", newText: \"432\u00f79=48, 0\" },\n  { oldText: \"306\u00f73=102, 0\", newText: \"170\u00f72=85, 0\" },\n  { oldText: \"770\u00f78=96, 2\", newText: \"418\u00f75=83, 3\" },\n  { oldText: \"761\u00f74=190, 1\", newText: \"926\u00f79=102, 8\" },\n  { oldText: \"316\u00f79=35, 1\", newText: \"600\u00f73=200, 0\" },\n  { oldText: \"118\u00f77=16, 6\", newText: \"388\u00f77=55, 3\" },\n  { oldText: \"536\u00f77=76, 4\", newText: \"920\u00f77=131, 3\" },\n  { oldText: \"886\u00f72=443, 0\", newText: \"462\u00f78=57, 6\" },\n  { oldText: \"452\u00f73=150, 2\", newText: \"118\u00f74=29, 2\" },\n  { oldText: \"920\u00f73=306, 2\", newText: \"471\u00f76=78, 3\" },\n  { oldText: \"179\u00f79=19, 8\", newText: \"950\u00f78=118, 6\" },\n  { oldText: \"215\u00f78=26, 7\", newText: \"366\u00f79=40, 6\" },\n  { oldText: \"162\u00f76=27, 0\", newText: \"649\u00f79=72, 1\" },\n  { oldText: \"764\u00f78=95, 4\", newText: \"444\u00f78=55, 4\" },\n  { oldText: \"247\u00f74=61, 3\", newText: \"751\u00f74=187, 3\" },\n  { oldText: \"317\u00f73=105, 2\", newText: \"273\u00f79=30, 3\" },\n  { oldText: \"579\u00f75=115, 4\", newText: \"204\u00f77=29, 1\" },\n  { oldText: \"516\u00f76=86, 0\", newText: \"211\u00f72=105, 1\" },\n  { oldText: \"724\u00f79=80, 4\", newText: \"570\u00f77=81, 3\" },\n  { oldText: \"496\u00f77=70, 6\", newText: \"928\u00f73=309, 1\" },\n  { oldText: \"621\u00f75=124, 1\", newText: \"825\u00f77=117, 6\" },\n  { oldText: \"914\u00f77=130, 4\", newText: \"329\u00f75=65, 4\" },\n  { oldText: \"823\u00f72=411, 1\", newText: \"394\u00f72=197, 0\" },\n  { oldText: \"210\u00f73=70, 0\", newText: \"793\u00f76=132, 1\" },\n];\n\nfor (const { oldText, newText } of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply literal text replacements for each unique cell/title text\n# using Word's Find/Replace (wdReplaceAll) on the whole document range.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-11-28 Thursday\"; New = \"2024-11-29 Friday\" }\n    @{ Old = \"598\u00f73=199, 1\"; New = \"947\u00f79=105, 2\" }\n    @{ Old = \"612\u00f74=153, 0\"; New = \"432\u00f79=48, 0\" }\n    @{ Old = \"306\u00f73=102, 0\"; New = \"170\u00f72=85, 0\" }\n    @{ Old = \"770\u00f78=96, 2\"; New = \"418\u00f75=83, 3\" }\n    @{ Old = \"761\u00f74=190, 1\"; New = \"926\u00f79=102, 8\" }\n    @{ Old = \"316\u00f79=35, 1\"; New = \"600\u00f73=200, 0\" }\n    @{ Old = \"118\u00f77=16, 6\"; New = \"388\u00f77=55, 3\" }\n    @{ Old = \"536\u00f77=76, 4\"; New = \"920\u00f77=131, 3\" }\n    @{ Old = \"886\u00f72=443, 0\"; New = \"462\u00f78=57, 6\" }\n    @{ Old = \"452\u00f73=150, 2\"; New = \"118\u00f74=29, 2\" }\n    @{ Old = \"920\u00f73=306, 2\"; New = \"471\u00f76=78, 3\" }\n    @{ Old = \"179\u00f79=19, 8\"; New = \"950\u00f78=118, 6\" }\n    @{ Old = \"215\u00f78=26, 7\"; New = \"366\u00f79=40, 6\" }\n    @{ Old = \"162\u00f76=27, 0\"; New = \"649\u00f79=72, 1\" }\n    @{ Old = \"764\u00f78=95, 4\"; New = \"444\u00f78=55, 4\" }\n    @{ Old = \"247\u00f74=61, 3\"; New = \"751\u00f74=187, 3\" }\n    @{ Old = \"317\u00f73=105, 2\"; New = \"273\u00f79=30, 3\" }\n    @{ Old = \"579\u00f75=115, 4\"; New = \"204\u00f77=29, 1\" }\n    @{ Old = \"516\u00f76=86, 0\"; New = \"211\u00f72=105, 1\" }\n    @{ Old = \"724\u00f79=80, 4\"; New = \"570\u00f77=81, 3\" }\n    @{ Old = \"496\u00f77=70, 6\"; New = \"928\u00f73=309, 1\" }\n    @{ Old = \"621\u00f75=124, 1\"; New = \"825\u00f77=117, 6\" }\n    @{ Old = \"914\u00f77=130, 4\"; New = \"329\u00f75=65, 4\" }\n    @{ Old = \"823\u00f72=411, 1\"; New = \"394\u00f72=197, 0\" }\n    @{ Old = \"210\u00f73=70, 0\"; New = \"793\u00f76=132, 1\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($r.Old, $False, $False, $False, $False, $False, $True, 1, $False, $r.New, 2) | Out-Null\n}\n"}
